$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns - copy formatting from the existing
# header cell (H1) so the new headers match the bold/centered/bordered style
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-33
$values = @(
    @(1, 2),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(7, 8),
    @(5, 6),
    @(9, 9),
    @(6, 7),
    @(7, 8),
    @(9, 9),
    @(5, 6),
    @(7, 8),
    @(5, 6),
    @(8, 9),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(6, 8),
    @(8, 8),
    @(3, 3),
    @(7, 7),
    @(2, 3),
    @(7, 7),
    @(9, 9),
    @(3, 4)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
